# Rename the three embedded picture parts in the headers/footers so the
# InlineShape "name" (wp:docPr/@name, mirrored onto pic:cNvPr/@name) swaps
# between the two logos that were already present in the package:
#
#   footer (Pearson logo, appears in both footers) : image1.png -> image2.png
#   header (BTec logo)                             : image2.jpg -> image1.jpg
#
# No pictures/parts are added or removed - this only touches the existing
# inline pictures' Name property.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers: both footer stories hold a copy of the Pearson logo --------
# Renaming straight through InlineShapes.Item(1).Name is unreliable for
# footer stories in this host, so round-trip through a floating Shape
# (ConvertToShape / ConvertToInlineShape) which reliably commits the
# rename and then put the picture back inline exactly as it was.
for ($f = 1; $f -le 2; $f++) {
    $footer = $sec.Footers.Item($f)
    if ($footer.Exists -and $footer.Range.InlineShapes.Count -ge 1) {
        $inlineShp = $footer.Range.InlineShapes.Item(1)
        if ($inlineShp.Name -ne "image2.png") {
            $floatShp = $inlineShp.ConvertToShape()
            $floatShp.Name = "image2.png"
            $floatShp.ConvertToInlineShape() | Out-Null
        }
    }
}

# --- Header: BTec logo (wdHeaderFooterFirstPage holds the picture) -------
for ($h = 1; $h -le 2; $h++) {
    $header = $sec.Headers.Item($h)
    if ($header.Exists -and $header.Range.InlineShapes.Count -ge 1) {
        $inlineShp = $header.Range.InlineShapes.Item(1)
        if ($inlineShp.Name -ne "image1.jpg") {
            $inlineShp.Name = "image1.jpg"
        }
    }
}
